$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value to a cell without Excel coercing
# numeric-looking strings (e.g. "1.002") into actual numbers, and without
# leaving a permanent style/number-format change behind on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '30.488.66'
Set-TextValue $ws.Range("E2") '  -0.97%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.111.61'
Set-TextValue $ws.Range("E3") '  -0.18%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  -0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '334.20'
Set-TextValue $ws.Range("E5") '  +0.43%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.001'
Set-TextValue $ws.Range("E6") '  +0.08%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5259'
Set-TextValue $ws.Range("E7") '  -1.35%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.4496'
Set-TextValue $ws.Range("E8") '  +2.23%  '

# Row 9
Set-TextValue $ws.Range("D9") '53.49'
Set-TextValue $ws.Range("E9") '  +12.63%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.09058'
Set-TextValue $ws.Range("E10") '  +0.57%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.177'
Set-TextValue $ws.Range("E11") '  -0.27%  '

# Row 12
Set-TextValue $ws.Range("D12") '24.49'
Set-TextValue $ws.Range("E12") '  -1.90%  '

# Row 13
Set-TextValue $ws.Range("D13") '2.101.52'
Set-TextValue $ws.Range("E13") '  -0.39%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.792'
Set-TextValue $ws.Range("E14") '  +0.38%  '

# Row 15
Set-TextValue $ws.Range("D15") '7.822'
Set-TextValue $ws.Range("E15") '  +0.27%  '

# Row 16
Set-TextValue $ws.Range("E16") '  +0.12%  '

# Row 17
Set-TextValue $ws.Range("E17") '  +0.03%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.00001131'
Set-TextValue $ws.Range("E18") '  -0.02%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06621'
Set-TextValue $ws.Range("E19") '  -0.90%  '

# Row 20
Set-TextValue $ws.Range("E20") '  +1.64%  '

# Row 21
Set-TextValue $ws.Range("E21") '  +0.06%  '

# Row 22
Set-TextValue $ws.Range("D22") '6.323'
Set-TextValue $ws.Range("E22") '  -0.19%  '

# Row 23
Set-TextValue $ws.Range("D23") '30.543.46'
Set-TextValue $ws.Range("E23") '  -1.00%  '

# Row 24
Set-TextValue $ws.Range("D24") '12.42'
Set-TextValue $ws.Range("E24") '  +0.63%  '

# Row 25
Set-TextValue $ws.Range("E25") '  +3.18%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.346.06'
Set-TextValue $ws.Range("E26") '  -0.52%  '

# Row 27
Set-TextValue $ws.Range("E27") '  -1.68%  '

# Row 28
Set-TextValue $ws.Range("E28") '  +0.33%  '

# Row 29
Set-TextValue $ws.Range("D29") '163.54'
Set-TextValue $ws.Range("E29") '  +0.30%  '

# Row 30
Set-TextValue $ws.Range("D30") '132.98'
Set-TextValue $ws.Range("E30") '  -0.42%  '

# Row 31
Set-TextValue $ws.Range("E31") '  +1.14%  '

# Row 32
Set-TextValue $ws.Range("E32") '  -0.75%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.669'
Set-TextValue $ws.Range("E33") '  +6.86%  '

# Row 34
Set-TextValue $ws.Range("D34") '6.171'
Set-TextValue $ws.Range("E34") '  -1.10%  '

# Row 35
Set-TextValue $ws.Range("D35") '3.925'
Set-TextValue $ws.Range("E35") '  -2.13%  '

# Row 36
Set-TextValue $ws.Range("D36") '10.62'
Set-TextValue $ws.Range("E36") '  +11.91%  '

# Row 37
Set-TextValue $ws.Range("E37") '  -0.62%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.06843'
Set-TextValue $ws.Range("E38") '  +1.27%  '

# Row 39
Set-TextValue $ws.Range("D39") '5.589'
Set-TextValue $ws.Range("E39") '  +1.10%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.2310'
Set-TextValue $ws.Range("E40") '  +0.41%  '

# Row 41
Set-TextValue $ws.Range("D41") '12.79'
Set-TextValue $ws.Range("E41") '  -0.84%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.6936'
Set-TextValue $ws.Range("E42") '  +1.42%  '

# Row 43
Set-TextValue $ws.Range("E43") '  -0.18%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.365'
Set-TextValue $ws.Range("E44") '  +6.00%  '

# Row 45
Set-TextValue $ws.Range("B45") 'EnergySwap'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '14.14'
Set-TextValue $ws.Range("E45") '  -0.15%  '

# Row 46
Set-TextValue $ws.Range("B46") 'Frax'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D46") '1.001'
Set-TextValue $ws.Range("E46") '  +0.09%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.6413'
Set-TextValue $ws.Range("E47") '  -0.52%  '

# Row 48
Set-TextValue $ws.Range("D48") '3.667'
Set-TextValue $ws.Range("E48") '  +0.21%  '

# Row 49
Set-TextValue $ws.Range("E49") '  -1.23%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.213'
Set-TextValue $ws.Range("E50") '  +2.34%  '

# Row 51
Set-TextValue $ws.Range("D51") '83.54'
Set-TextValue $ws.Range("E51") '  +0.63%  '
